{"js": "// Replace each three-digit-by-one-digit multiplication prompt in the\n// practice table with its updated problem, e.g. \"497\u00d74=\" -> \"693\u00d76=\".\n// Every old value in the document is unique, so a plain matchCase\n// search-and-replace on the body is unambiguous.\nconst replacements = [\n  [\"497\u00d74=\", \"693\u00d76=\"],\n  [\"678\u00d74=\", \"714\u00d72=\"],\n  [\"741\u00d76=\", \"832\u00d74=\"],\n  [\"289\u00d77=\", \"448\u00d79=\"],\n  [\"729\u00d76=\", \"313\u00d77=\"],\n  [\"366\u00d77=\", \"274\u00d72=\"],\n  [\"535\u00d75=\", \"620\u00d78=\"],\n  [\"929\u00d73=\", \"781\u00d72=\"],\n  [\"734\u00d73=\", \"529\u00d78=\"],\n  [\"716\u00d74=\", \"955\u00d79=\"],\n  [\"294\u00d72=\", \"461\u00d76=\"],\n  [\"960\u00d79=\", \"673\u00d74=\"],\n  [\"169\u00d79=\", \"302\u00d76=\"],\n  [\"503\u00d79=\", \"791\u00d77=\"],\n  [\"673\u00d72=\", \"843\u00d77=\"],\n  [\"968\u00d78=\", \"674\u00d78=\"],\n  [\"365\u00d76=\", \"211\u00d79=\"],\n  [\"264\u00d78=\", \"526\u00d78=\"],\n  [\"617\u00d75=\", \"305\u00d72=\"],\n  [\"824\u00d75=\", \"439\u00d73=\"],\n  [\"659\u00d76=\", \"117\u00d74=\"],\n  [\"154\u00d72=\", \"771\u00d72=\"],\n  [\"903\u00d78=\", \"468\u00d72=\"],\n  [\"397\u00d74=\", \"525\u00d75=\"],\n  [\"809\u00d75=\", \"431\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication prompt in the\n# practice table with its updated problem, e.g. \"497\u00d74=\" -> \"693\u00d76=\".\n# Every old value in the document is unique, so Find/Execute with\n# MatchCase=True and no wildcards resolves each pair unambiguously.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"497\u00d74=\", \"693\u00d76=\"),\n  @(\"678\u00d74=\", \"714\u00d72=\"),\n  @(\"741\u00d76=\", \"832\u00d74=\"),\n  @(\"289\u00d77=\", \"448\u00d79=\"),\n  @(\"729\u00d76=\", \"313\u00d77=\"),\n  @(\"366\u00d77=\", \"274\u00d72=\"),\n  @(\"535\u00d75=\", \"620\u00d78=\"),\n  @(\"929\u00d73=\", \"781\u00d72=\"),\n  @(\"734\u00d73=\", \"529\u00d78=\"),\n  @(\"716\u00d74=\", \"955\u00d79=\"),\n  @(\"294\u00d72=\", \"461\u00d76=\"),\n  @(\"960\u00d79=\", \"673\u00d74=\"),\n  @(\"169\u00d79=\", \"302\u00d76=\"),\n  @(\"503\u00d79=\", \"791\u00d77=\"),\n  @(\"673\u00d72=\", \"843\u00d77=\"),\n  @(\"968\u00d78=\", \"674\u00d78=\"),\n  @(\"365\u00d76=\", \"211\u00d79=\"),\n  @(\"264\u00d78=\", \"526\u00d78=\"),\n  @(\"617\u00d75=\", \"305\u00d72=\"),\n  @(\"824\u00d75=\", \"439\u00d73=\"),\n  @(\"659\u00d76=\", \"117\u00d74=\"),\n  @(\"154\u00d72=\", \"771\u00d72=\"),\n  @(\"903\u00d78=\", \"468\u00d72=\"),\n  @(\"397\u00d74=\", \"525\u00d75=\"),\n  @(\"809\u00d75=\", \"431\u00d77=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $range.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
